$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1501.3472
$ws.Range("J17").Value = 1501.3472
$ws.Range("L17").Value = 4504.0416
$ws.Range("N17").Value = -4840.0416

$ws.Range("H43").Value = 3057.4285
$ws.Range("I43").Value = 3057.4285
$ws.Range("K43").Value = 3057.4285
$ws.Range("M43").Value = -2988.4285

$ws.Range("H51").Value = 7374.125
$ws.Range("I51").Value = 8332.5
$ws.Range("J51").Value = 4499
$ws.Range("K51").Value = 8332.5
$ws.Range("L51").Value = 4499
$ws.Range("M51").Value = -7848.5
$ws.Range("N51").Value = -5467

$ws.Range("H62").Value = 3570.4285
$ws.Range("I62").Value = 3448.5
$ws.Range("J62").Value = 3733
$ws.Range("K62").Value = 3448.5
$ws.Range("L62").Value = 3733
$ws.Range("M62").Value = -2824.5
$ws.Range("N62").Value = -4981

$ws.Range("H65").Value = 3570.4285
$ws.Range("I65").Value = 3448.5
$ws.Range("J65").Value = 3733
$ws.Range("K65").Value = 17242.5
$ws.Range("L65").Value = 18665
$ws.Range("M65").Value = -14122.5
$ws.Range("N65").Value = -24905

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H80").Value = 6997.3335
$ws.Range("I80").Value = 5362
$ws.Range("J80").Value = 8632.666999999999
$ws.Range("K80").Value = 16086
$ws.Range("L80").Value = 25898.001
$ws.Range("M80").Value = -15088
$ws.Range("N80").Value = -27894.001

$ws.Range("H83").Value = 6997.3335
$ws.Range("I83").Value = 5362
$ws.Range("J83").Value = 8632.666999999999
$ws.Range("K83").Value = 48258
$ws.Range("L83").Value = 77694.003
$ws.Range("M83").Value = -43266
$ws.Range("N83").Value = -87678.003

$ws.Range("H92").Value = 91483.63
$ws.Range("I92").Value = 91483.63
$ws.Range("K92").Value = 91483.63
$ws.Range("M92").Value = -90235.63

$ws.Range("H106").Value = 26777.889
$ws.Range("I106").Value = 25286.285
$ws.Range("K106").Value = 25286.285
$ws.Range("M106").Value = -24655.285

$ws.Range("H116").Value = 5000
$ws.Range("I116").Value = 5000
$ws.Range("K116").Value = 5000
$ws.Range("M116").Value = -1558

$ws.Range("H137").Value = 2230.4546
$ws.Range("I137").Value = 1751.3889
$ws.Range("J137").Value = 4386.25
$ws.Range("K137").Value = 5254.1667
$ws.Range("L137").Value = 13158.75
$ws.Range("M137").Value = -2704.1667
$ws.Range("N137").Value = -18258.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 3999.6667
$ws.Range("I31").Value = 3999.6667
$ws.Range("K31").Value = 3999.6667
$ws.Range("M31").Value = -3705.6667

$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 28481
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 28481
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -29107

$ws.Range("H45").Value = 2139.4666
$ws.Range("I45").Value = 2183.5
$ws.Range("J45").Value = 1963.3334
$ws.Range("K45").Value = 2183.5
$ws.Range("L45").Value = 1963.3334
$ws.Range("M45").Value = -1806.5
$ws.Range("N45").Value = -2717.3334

$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 257.1111
$ws.Range("I22").Value = 282.42856
$ws.Range("K22").Value = 282.42856
$ws.Range("M22").Value = -109.42856

$ws.Range("H59").Value = 99999
$ws.Range("J59").Value = 99999
$ws.Range("L59").Value = 99999
$ws.Range("N59").Value = -101693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5000446.5
$ws.Range("I22").Value = 524.2
$ws.Range("J22").Value = 13333650
$ws.Range("K22").Value = 524.2
$ws.Range("L22").Value = 13333650
$ws.Range("M22").Value = -174.2
$ws.Range("N22").Value = -13334350

$ws.Range("H86").Value = 7998.4287
$ws.Range("I86").Value = 7494.5
$ws.Range("K86").Value = 7494.5
$ws.Range("M86").Value = -6371.5

$ws.Range("H89").Value = 7998.4287
$ws.Range("I89").Value = 7494.5
$ws.Range("K89").Value = 37472.5
$ws.Range("M89").Value = -31856.5

$ws.Range("H107").Value = 1331.9231
$ws.Range("I107").Value = 610.1429000000001
$ws.Range("J107").Value = 2174
$ws.Range("K107").Value = 610.1429000000001
$ws.Range("L107").Value = 2174
$ws.Range("M107").Value = 1309.8571
$ws.Range("N107").Value = -6014

$ws.Range("H121").Value = 70326
$ws.Range("J121").Value = 70326
$ws.Range("L121").Value = 70326
$ws.Range("N121").Value = -72946

$ws.Range("H134").Value = 3416.5557
$ws.Range("I134").Value = 3416.5557
$ws.Range("K134").Value = 10249.6671
$ws.Range("M134").Value = -7714.667099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1144.8889
$ws.Range("I5").Value = 1046.7142
$ws.Range("J5").Value = 1488.5
$ws.Range("K5").Value = 3140.1426
$ws.Range("L5").Value = 4465.5
$ws.Range("M5").Value = -3028.1426
$ws.Range("N5").Value = -4689.5

$ws.Range("H129").Value = 2813.6667
$ws.Range("J129").Value = 2813.6667
$ws.Range("L129").Value = 8441.000100000001
$ws.Range("N129").Value = -18441.0001

$ws.Range("H131").Value = 2287.5
$ws.Range("I131").Value = 1681
$ws.Range("J131").Value = 3136.6
$ws.Range("K131").Value = 5043
$ws.Range("L131").Value = 9409.799999999999
$ws.Range("M131").Value = -3
$ws.Range("N131").Value = -19489.8

$ws.Range("H135").Value = 1144.8889
$ws.Range("I135").Value = 1046.7142
$ws.Range("J135").Value = 1488.5
$ws.Range("K135").Value = 9420.427799999999
$ws.Range("L135").Value = 13396.5
$ws.Range("M135").Value = -6885.427799999999
$ws.Range("N135").Value = -18466.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 5000
$ws.Range("I5").Value = 5000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 5000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -4888
$ws.Range("N5").ClearContents()

$ws.Range("H57").Value = 2989.6
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H70").Value = 16929.666
$ws.Range("J70").Value = 6740
$ws.Range("L70").Value = 6740
$ws.Range("N70").Value = -7280

$ws.Range("H73").Value = 16929.666
$ws.Range("J73").Value = 6740
$ws.Range("L73").Value = 6740
$ws.Range("N73").Value = -8612

$ws.Range("H126").Value = 2937.5
$ws.Range("J126").Value = 2750
$ws.Range("L126").Value = 8250
$ws.Range("N126").Value = -13190

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1731.2667
$ws.Range("I61").Value = 1755
$ws.Range("K61").Value = 1755
$ws.Range("M61").Value = -1553

$ws.Range("H94").Value = 49750
$ws.Range("J94").Value = 49750
$ws.Range("L94").Value = 49750
$ws.Range("N94").Value = -51102

$ws.Range("H113").Value = 1731.2667
$ws.Range("I113").Value = 1755
$ws.Range("K113").Value = 1755
$ws.Range("M113").Value = 415

$ws.Range("H122").Value = 4460.1113
$ws.Range("I122").Value = 3907.0833
$ws.Range("K122").Value = 11721.2499
$ws.Range("M122").Value = -9271.249899999999

$ws.Range("H136").Value = 3148
$ws.Range("I136").Value = 1555.75
$ws.Range("K136").Value = 4667.25
$ws.Range("M136").Value = -2117.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 44996.25
$ws.Range("I15").Value = 44990
$ws.Range("K15").Value = 44990
$ws.Range("M15").Value = -44702

$ws.Range("H46").Value = 50013696
$ws.Range("I46").Value = 27390
$ws.Range("K46").Value = 27390
$ws.Range("M46").Value = -27159

$ws.Range("H134").Value = 50013696
$ws.Range("I134").Value = 27390
$ws.Range("K134").Value = 82170
$ws.Range("M134").Value = -79635

$ws.Range("H136").Value = 7858.28
$ws.Range("I136").Value = 5194.625
$ws.Range("K136").Value = 15583.875
$ws.Range("M136").Value = -13033.875
